# Documentation: consolidated installation instructions. (#40)
#
# 1) The cached "today" text of every Date placeholder (slide master +
#    every custom layout) gets bumped from 9/11/2019 to 9/16/2019 -
#    PowerPoint recaches the datetimeFigureOut field's visible text
#    whenever the deck is saved on a later day.
# 2) The "6. sync" instructions textbox on the single content slide grows
#    by one line ("git push origin <i>feature</i>") and is nudged up to
#    keep its connector anchors lined up.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, $newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "9/16/2019"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "9/16/2019"
}

# --- "6. sync" textbox on slide 1 -----------------------------------

$s = $p.Slides.Item(1)

$syncBox = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 96") {
        $syncBox = $shp
    }
}

$tr = $syncBox.TextFrame.TextRange
$tr.InsertAfter("`rgit push origin feature") | Out-Null

# The new paragraph inherited the italic formatting of the preceding
# "message" run; split it so only "feature" stays italic.
$full = $syncBox.TextFrame.TextRange
$paraCount = $full.Paragraphs().Count
$lastPara = $full.Paragraphs($paraCount)
$featureLen = 7
$prefixLen = $lastPara.Length - $featureLen
$prefixRange = $full.Characters($lastPara.Start, $prefixLen)
$prefixRange.Font.Italic = $false

# Box grew (autosize) to fit the extra line; also nudge it up so it
# stays lined up with its connector, matching the authored layout.
$syncBox.Top = 195.4487762451172
